$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035200265770164
$ws.Range("D2").Value = 1.044635980063504
$ws.Range("E2").Value = 1.052638256155214
$ws.Range("F2").Value = 1.058119108826274
$ws.Range("I2").Value = 1.037241743938605
$ws.Range("J2").Value = 1.040315340478905
$ws.Range("K2").Value = 1.047406473663439
$ws.Range("L2").Value = 1.055386435026599
$ws.Range("M2").Value = 1.060852215044985
$ws.Range("N2").Value = 1.041792707588223

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036155896752796
$ws.Range("D3").Value = 1.04537756060644
$ws.Range("E3").Value = 1.053550349791475
$ws.Range("F3").Value = 1.059010151720544
$ws.Range("I3").Value = 1.037407022813624
$ws.Range("J3").Value = 1.040914389777309
$ws.Range("K3").Value = 1.047959514080242
$ws.Range("L3").Value = 1.05611115965238
$ws.Range("M3").Value = 1.061557032709588
$ws.Range("N3").Value = 1.042392607605342

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036774746733029
$ws.Range("D4").Value = 1.045857702897246
$ws.Range("E4").Value = 1.054141679311682
$ws.Range("F4").Value = 1.059587511591875
$ws.Range("I4").Value = 1.037512773598052
$ws.Range("J4").Value = 1.041301895839258
$ws.Range("K4").Value = 1.04831699144894
$ws.Range("L4").Value = 1.056580614999699
$ws.Range("M4").Value = 1.062013253780765
$ws.Range("N4").Value = 1.042780663970347

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037035028294216
$ws.Range("D5").Value = 1.046059622740226
$ws.Range("E5").Value = 1.054390546609109
$ws.Range("F5").Value = 1.059830422635389
$ws.Range("I5").Value = 1.037556944423396
$ws.Range("J5").Value = 1.041464773957415
$ws.Range("K5").Value = 1.048467183661131
$ws.Range("L5").Value = 1.056778094880707
$ws.Range("M5").Value = 1.062205085423169
$ws.Range("N5").Value = 1.042943773394113

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037078737548626
$ws.Range("D6").Value = 1.04609352991789
$ws.Range("E6").Value = 1.054432348451689
$ws.Range("F6").Value = 1.059871219505223
$ws.Range("I6").Value = 1.037564344064978
$ws.Range("J6").Value = 1.041492120156253
$ws.Range("K6").Value = 1.048492396210842
$ws.Range("L6").Value = 1.056811259672107
$ws.Range("M6").Value = 1.062237296892283
$ws.Range("N6").Value = 1.042971158427691

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036778224169329
$ws.Range("D7").Value = 1.045860400694572
$ws.Range("E7").Value = 1.05414500362035
$ws.Range("F7").Value = 1.059590756640549
$ws.Range("I7").Value = 1.037513364938413
$ws.Range("J7").Value = 1.041304072341497
$ws.Range("K7").Value = 1.048318998684622
$ws.Range("L7").Value = 1.056583253262238
$ws.Range("M7").Value = 1.062015816903491
$ws.Range("N7").Value = 1.042782843563469

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035523123114653
$ws.Range("D8").Value = 1.044886539758375
$ws.Range("E8").Value = 1.052946264708122
$ws.Range("F8").Value = 1.05842007532138
$ws.Range("I8").Value = 1.037297847953807
$ws.Range("J8").Value = 1.040517816029301
$ws.Range("K8").Value = 1.047593453625484
$ws.Range("L8").Value = 1.055631252742315
$ws.Range("M8").Value = 1.061090378226454
$ws.Range("N8").Value = 1.041995470677124

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.0333152914302
$ws.Range("D9").Value = 1.043172756994837
$ws.Range("E9").Value = 1.050842760990109
$ws.Range("F9").Value = 1.056363344364621
$ws.Range("I9").Value = 1.036908946241351
$ws.Range("J9").Value = 1.039131464089577
$ws.Range("K9").Value = 1.046312112492944
$ws.Range("L9").Value = 1.053957668042603
$ws.Range("M9").Value = 1.05946089450463
$ws.Range("N9").Value = 1.040607149958629

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031846025414848
$ws.Range("D10").Value = 1.042031855422512
$ws.Range("E10").Value = 1.049446448861814
$ws.Range("F10").Value = 1.054996420117855
$ws.Range("I10").Value = 1.036643570364821
$ws.Range("J10").Value = 1.038206702851221
$ws.Range("K10").Value = 1.045456040418369
$ws.Range("L10").Value = 1.052844682897175
$ws.Range("M10").Value = 1.058375492004066
$ws.Range("N10").Value = 1.039681075453243

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.031210450157135
$ws.Range("D11").Value = 1.041538234332009
$ws.Range("E11").Value = 1.048843277561924
$ws.Range("F11").Value = 1.054405548570408
$ws.Range("I11").Value = 1.036527217704487
$ws.Range("J11").Value = 1.037806157888204
$ws.Range("K11").Value = 1.045084927428797
$ws.Range("L11").Value = 1.052363412957898
$ws.Range("M11").Value = 1.057905734474637
$ws.Range("N11").Value = 1.039279961670435

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030974464206158
$ws.Range("D12").Value = 1.041354942493253
$ws.Range("E12").Value = 1.048619450869396
$ws.Range("F12").Value = 1.054186226633747
$ws.Range("I12").Value = 1.036483782846559
$ws.Range("J12").Value = 1.037657360955644
$ws.Range("K12").Value = 1.044947016251715
$ws.Range("L12").Value = 1.052184748258511
$ws.Range("M12").Value = 1.057731281230304
$ws.Range("N12").Value = 1.039130953429164

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.031025079681029
$ws.Range("D13").Value = 1.041394256431515
$ws.Range("E13").Value = 1.048667452566701
$ws.Range("F13").Value = 1.054233264938282
$ws.Range("I13").Value = 1.036493109550687
$ws.Range("J13").Value = 1.037689279143731
$ws.Range("K13").Value = 1.044976601514712
$ws.Range("L13").Value = 1.052223067880628
$ws.Range("M13").Value = 1.057768700411919
$ws.Range("N13").Value = 1.039162916944739

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.031190941541896
$ws.Range("D14").Value = 1.04152308211873
$ws.Range("E14").Value = 1.048824771522312
$ws.Range("F14").Value = 1.05438741619932
$ws.Range("I14").Value = 1.036523631775016
$ws.Range("J14").Value = 1.037793858622694
$ws.Range("K14").Value = 1.045073528932974
$ws.Range("L14").Value = 1.052348642423171
$ws.Range("M14").Value = 1.057891313381346
$ws.Range("N14").Value = 1.039267644938558

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031293147122458
$ws.Range("D15").Value = 1.041602464002913
$ws.Range("E15").Value = 1.048921729852411
$ws.Range("F15").Value = 1.054482414337969
$ws.Range("I15").Value = 1.036542408878569
$ws.Range("J15").Value = 1.037858291305876
$ws.Range("K15").Value = 1.045133240758893
$ws.Range("L15").Value = 1.052426026382096
$ws.Range("M15").Value = 1.057966864034983
$ws.Range("N15").Value = 1.039332169123541

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031888219726625
$ws.Range("D16").Value = 1.042064623908207
$ws.Range("E16").Value = 1.049486509842446
$ws.Range("F16").Value = 1.055035655841079
$ws.Range("I16").Value = 1.036651261946292
$ws.Range("J16").Value = 1.038233283309167
$ws.Range("K16").Value = 1.045480661033162
$ws.Range("L16").Value = 1.052876637207954
$ws.Range("M16").Value = 1.058406673191801
$ws.Range("N16").Value = 1.039707693658488

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032261661416779
$ws.Range("D17").Value = 1.042354631831504
$ws.Range("E17").Value = 1.049841168268011
$ws.Range("F17").Value = 1.05538296258212
$ws.Range("I17").Value = 1.036719156591808
$ws.Range("J17").Value = 1.038468475141944
$ws.Range("K17").Value = 1.045698474698389
$ws.Range("L17").Value = 1.053159471036333
$ws.Range("M17").Value = 1.058682616008686
$ws.Range("N17").Value = 1.039943219490644

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032479544038832
$ws.Range("D18").Value = 1.04252382674037
$ws.Range("E18").Value = 1.050048173469902
$ws.Range("F18").Value = 1.055585638627061
$ws.Range("I18").Value = 1.036758619111505
$ws.Range("J18").Value = 1.038605647185004
$ws.Range("K18").Value = 1.045825480346328
$ws.Range("L18").Value = 1.053324506818481
$ws.Range("M18").Value = 1.058843590789959
$ws.Range("N18").Value = 1.040080586333739

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032553846601686
$ws.Range("D19").Value = 1.042581524267778
$ws.Range("E19").Value = 1.050118780416305
$ws.Range("F19").Value = 1.055654762481364
$ws.Range("I19").Value = 1.03677205117115
$ws.Range("J19").Value = 1.038652417349809
$ws.Range("K19").Value = 1.045868778940338
$ws.Range("L19").Value = 1.053380790504706
$ws.Range("M19").Value = 1.058898482725525
$ws.Range("N19").Value = 1.040127422917543

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032221588403726
$ws.Range("D20").Value = 1.042323512765291
$ws.Range("E20").Value = 1.049803102393213
$ws.Range("F20").Value = 1.055345689718904
$ws.Range("I20").Value = 1.036711886541356
$ws.Range("J20").Value = 1.038443242453938
$ws.Range("K20").Value = 1.045675109625907
$ws.Range("L20").Value = 1.053129119056406
$ws.Range("M20").Value = 1.058653007671004
$ws.Range("N20").Value = 1.039917950969327

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03114209672854
$ws.Range("D21").Value = 1.041485144469289
$ws.Range("E21").Value = 1.048778438981622
$ws.Range("F21").Value = 1.054342018230304
$ws.Range("I21").Value = 1.036514649707446
$ws.Range("J21").Value = 1.037763063026329
$ws.Range("K21").Value = 1.045044987962932
$ws.Range("L21").Value = 1.052311661069256
$ws.Range("M21").Value = 1.057855205924813
$ws.Range("N21").Value = 1.039236805608913

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.030463927528432
$ws.Range("D22").Value = 1.040958382738063
$ws.Range("E22").Value = 1.048135454350665
$ws.Range("F22").Value = 1.053711862229393
$ws.Range("I22").Value = 1.036389387898549
$ws.Range("J22").Value = 1.037335311062963
$ws.Range("K22").Value = 1.044648439805744
$ws.Range("L22").Value = 1.051798274192617
$ws.Range("M22").Value = 1.057353802985262
$ws.Range("N22").Value = 1.038808446188695

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030823385426193
$ws.Range("D23").Value = 1.041237595106864
$ws.Range("E23").Value = 1.048476192558606
$ws.Range("F23").Value = 1.054045834796071
$ws.Range("I23").Value = 1.036455909977658
$ws.Range("J23").Value = 1.037562079280875
$ws.Range("K23").Value = 1.044858691801782
$ws.Range("L23").Value = 1.05207037476768
$ws.Range("M23").Value = 1.057619586169497
$ws.Range("N23").Value = 1.039035536443488

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.032239695486808
$ws.Range("D24").Value = 1.04233757401394
$ws.Range("E24").Value = 1.049820302294356
$ws.Range("F24").Value = 1.055362531419059
$ws.Range("I24").Value = 1.036715171994613
$ws.Range("J24").Value = 1.038454644054738
$ws.Range("K24").Value = 1.045685667424369
$ws.Range("L24").Value = 1.053142833613276
$ws.Range("M24").Value = 1.058666386336454
$ws.Range("N24").Value = 1.039929368761708

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033885610826464
$ws.Range("D25").Value = 1.043615531498346
$ws.Range("E25").Value = 1.05138551190055
$ws.Range("F25").Value = 1.056894319596044
$ws.Range("I25").Value = 1.037010565638022
$ws.Range("J25").Value = 1.039489966285671
$ws.Range("K25").Value = 1.046643699672771
$ws.Range("L25").Value = 1.054389852655443
$ws.Range("M25").Value = 1.059881998413344
$ws.Range("N25").Value = 1.040966161268962
